# This edit permutes the "species observation" data among the existing
# data rows (7-13, 15-21) of the Artfynd sheet. Row 14 is untouched.
# Every row keeps its fixed/common columns (C, P, S, T, U, V, W, Y, Z, AA,
# AB, AD, AE, AG, AT, AW, AX, AY) and only the per-observation fields
# (A, B, D, E, F, G, H, Q, R) plus the optional K/L/M/N "empty marker"
# cells and the AC "ringhack äldre" comment move between row positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(7, 8, 9, 10, 11, 12, 13, 15, 16, 17, 18, 19, 20, 21)

# Snapshot the current ("before") content of every mobile field for each
# row so that writes to one row never clobber data still needed for
# another row later in the script.
$snapshot = @{}
foreach ($r in $rows) {
    $row = @{}
    $row["A"] = $ws.Range("A$r").Value2
    $row["B"] = $ws.Range("B$r").Value2
    $row["D"] = $ws.Range("D$r").Value2
    $row["E"] = $ws.Range("E$r").Value2
    $row["F"] = $ws.Range("F$r").Value2
    $row["G"] = $ws.Range("G$r").Value2
    $row["H"] = $ws.Range("H$r").Value2
    $row["Q"] = $ws.Range("Q$r").Value2
    $row["R"] = $ws.Range("R$r").Value2
    $row["AC"] = $ws.Range("AC$r").Value2
    $snapshot[$r] = $row
}

# Rows that originally carried the blank K/L/M/N marker cells together
# with the AC "ringhack äldre" comment (the three "Tretåig hackspett"
# observations).
$specialSource = @(8, 16, 21)

# Mapping of target row -> source row: the target row's mobile data
# (after the edit) is the data that the source row held beforehand.
$mapping = @{
    7  = 12
    8  = 13
    9  = 18
    10 = 16
    11 = 17
    12 = 19
    13 = 8
    15 = 21
    16 = 15
    17 = 11
    18 = 7
    19 = 10
    20 = 9
    21 = 20
}

# First, physically copy the blank K/L/M/N "marker" cells from each
# special source row to whichever target row will receive its data.
# This has to happen before any of the source/target cells are
# overwritten below, and Copy() (unlike assigning Value = "") actually
# materializes an empty cell rather than removing it.
foreach ($t in $rows) {
    $s = $mapping[$t]
    if ($specialSource -contains $s) {
        $ws.Range("K$s").Copy($ws.Range("K$t"))
        $ws.Range("L$s").Copy($ws.Range("L$t"))
        $ws.Range("M$s").Copy($ws.Range("M$t"))
        $ws.Range("N$s").Copy($ws.Range("N$t"))
    }
}

foreach ($t in $rows) {
    $s = $mapping[$t]
    $data = $snapshot[$s]

    $ws.Range("A$t").Value = $data["A"]
    $ws.Range("B$t").Value = $data["B"]
    $ws.Range("D$t").Value = $data["D"]
    $ws.Range("E$t").Value = $data["E"]
    $ws.Range("F$t").Value = $data["F"]
    $ws.Range("G$t").Value = $data["G"]
    $ws.Range("H$t").Value = $data["H"]
    $ws.Range("Q$t").Value = $data["Q"]
    $ws.Range("R$t").Value = $data["R"]

    if ($specialSource -contains $s) {
        # Target row should now carry the AC comment (the K/L/M/N
        # markers were already copied into place above).
        $ws.Range("AC$t").Value = "ringhack äldre"
    }
    else {
        # Target row should not have the K/L/M/N markers nor the AC
        # comment - clear them in case the row held them beforehand.
        $ws.Range("K$t").ClearContents()
        $ws.Range("L$t").ClearContents()
        $ws.Range("M$t").ClearContents()
        $ws.Range("N$t").ClearContents()
        $ws.Range("AC$t").ClearContents()
    }
}
